# Data-message workbook update
# - Rows 2,4,8 (column A) and rows 3,5,9 (column A) lost the surrounding
#   "*...*" emphasis markers around the two prayer texts, matching the
#   already-unmarked text that rows 6/7 used.
# - Row 8/9 (column B) time-window moved from 17:45-17:54 to 21:00-21:09.
# - Selection moved to A13 with the view scrolled so row 7 is at the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull the already "plain" (no-asterisk) prayer texts from rows 6/7 so the
# exact wording/line breaks are reused verbatim for the other occurrences.
$pokayanie = $ws.Range("A6").Value2
$otcheNash = $ws.Range("A7").Value2

$ws.Range("A2").Value = $pokayanie
$ws.Range("A3").Value = $otcheNash
$ws.Range("A4").Value = $pokayanie
$ws.Range("A5").Value = $otcheNash
$ws.Range("A8").Value = $pokayanie
$ws.Range("A9").Value = $otcheNash

# Row 4 has a fixed custom height (not auto-expanded like rows 2/6/8, which
# already sit at Excel's 409.5pt cap); re-pin it since the longer text would
# otherwise trigger an autofit growth that the source file does not show.
$ws.Rows(4).RowHeight = 79.5

# Update the time window on row 8/9 from 17:45-17:54 to 21:00-21:09.
$ws.Range("B8").Value = "21:00 - 21:04"
$ws.Range("B9").Value = "21:05 - 21:09"

# Reflect the updated selection/scroll position from the saved view.
$ws.Activate()
$ws.Range("A13").Select()
